$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D:E (numeric/percent-looking strings) to avoid Excel auto-typing them as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "31.246.47"
$ws.Range("E2").Value = "  +2.16%  "

$ws.Range("D3").Value = "1.996.95"
$ws.Range("E3").Value = "  +6.05%  "

$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "0.7760"
$ws.Range("E5").Value = "  +63.72%  "

$ws.Range("D6").Value = "254.04"
$ws.Range("E6").Value = "  +3.20%  "

$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "0.3483"
$ws.Range("E8").Value = "  +20.42%  "

$ws.Range("D9").Value = "27.96"
$ws.Range("E9").Value = "  +25.36%  "

$ws.Range("D10").Value = "0.07077"
$ws.Range("E10").Value = "  +8.37%  "

$ws.Range("D11").Value = "0.8405"
$ws.Range("E11").Value = "  +10.41%  "

$ws.Range("D12").Value = "0.08205"
$ws.Range("E12").Value = "  +4.88%  "

$ws.Range("D13").Value = "100.92"
$ws.Range("E13").Value = "  +1.42%  "

$ws.Range("D14").Value = "1.994.82"
$ws.Range("E14").Value = "  +5.97%  "

$ws.Range("D15").Value = "5.629"
$ws.Range("E15").Value = "  +7.65%  "

$ws.Range("D16").Value = "15.24"
$ws.Range("E16").Value = "  +15.81%  "

$ws.Range("D17").Value = "272.25"
$ws.Range("E17").Value = "  -3.98%  "

$ws.Range("D18").Value = "31.240.54"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("D19").Value = "5.962"
$ws.Range("E19").Value = "  +11.47%  "

$ws.Range("D20").Value = "0.000008010"
$ws.Range("E20").Value = "  +6.56%  "

$ws.Range("D21").Value = "2.256.73"
$ws.Range("E21").Value = "  +6.07%  "

$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "7.095"
$ws.Range("E24").Value = "  +10.36%  "

$ws.Range("D25").Value = "10.01"
$ws.Range("E25").Value = "  +9.23%  "

$ws.Range("D26").Value = "164.59"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").Value = "0.1422"
$ws.Range("E27").Value = "  +45.91%  "

$ws.Range("D28").Value = "2.428"
$ws.Range("E28").Value = "  +27.50%  "

$ws.Range("D29").Value = "19.88"
$ws.Range("E29").Value = "  +4.60%  "

$ws.Range("D30").Value = "1.594"
$ws.Range("E30").Value = "  +6.15%  "

$ws.Range("D31").Value = "1.363"
$ws.Range("E31").Value = "  +2.70%  "

$ws.Range("D32").Value = "4.611"
$ws.Range("E32").Value = "  +8.58%  "

$ws.Range("D33").Value = "4.447"
$ws.Range("E33").Value = "  +6.38%  "

$ws.Range("D34").Value = "0.05325"
$ws.Range("E34").Value = "  +10.03%  "

$ws.Range("D35").Value = "1.246"
$ws.Range("E35").Value = "  +10.23%  "

$ws.Range("D36").Value = "0.7943"
$ws.Range("E36").Value = "  +13.83%  "

$ws.Range("D37").Value = "2.768"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "0.9990"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("D39").Value = "0.02003"
$ws.Range("E39").Value = "  +5.20%  "

$ws.Range("D40").Value = "2.917"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").Value = "83.46"
$ws.Range("E41").Value = "  +10.74%  "

$ws.Range("D42").Value = "6.777"
$ws.Range("E42").Value = "  +7.53%  "

$ws.Range("D43").Value = "0.4677"
$ws.Range("E43").Value = "  +10.18%  "

$ws.Range("D44").Value = "2.136"
$ws.Range("E44").Value = "  +8.30%  "

$ws.Range("D45").Value = "0.8548"
$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("D46").Value = "105.00"
$ws.Range("E46").Value = "  +3.58%  "

$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").Value = "7.767"
$ws.Range("E48").Value = "  +10.80%  "

$ws.Range("D49").Value = "10.04"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "37.59"
$ws.Range("E50").Value = "  +6.68%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.4334"
$ws.Range("E51").Value = "  +9.67%  "

# Restore default (unstyled) formatting on the touched numeric/percent columns so number format
# metadata does not leak into the saved styles.
$ws.Range("D2:E51").Style = "Normal"
